$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '244.04'

# Row 3
Set-TextValue $ws.Range('D3') '25.00'

# Row 4
$ws.Range('B4').Value = 'LEO'
$ws.Range('C4').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D4') '3.500'
$ws.Range('E4').Value = '3LEOLEO'

# Row 5
$ws.Range('B5').Value = 'HuobiToken'
$ws.Range('C5').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D5') '5.174'
$ws.Range('E5').Value = '4HuobiTokenHT'

# Row 6
$ws.Range('B6').Value = 'Cronos'
$ws.Range('C6').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D6') '0.05669'
$ws.Range('E6').Value = '5CronosCRO'

# Row 7
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws.Range('D7') '6.517'
$ws.Range('E7').Value = '6KuCoinTokenKCS'

# Row 8
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D8') '2.972'
$ws.Range('E8').Value = '7GateTokenGT'

# Row 9
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D9') '0.8093'
$ws.Range('E9').Value = '8MXTokenMX'

# Row 10
$ws.Range('B10').Value = 'FTXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D10') '0.8347'
$ws.Range('E10').Value = '9FTXTokenFTT'

# Row 11
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D11') '0.1332'
$ws.Range('E11').Value = '10WazirXWRX'

# Row 12
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D12') '0.06942'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D13') '0.02826'
$ws.Range('E13').Value = '12BitrueCoinBTR'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D14') '0.09384'
$ws.Range('E14').Value = '13BitMartTokenBMX'

# Row 15
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D15') '0.001509'
$ws.Range('E15').Value = '14BitForexTokenBF'

# Row 16
Set-TextValue $ws.Range('D16') '0.006081'

# Row 17
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D17') '2.091'
$ws.Range('E17').Value = '16BTSETokenBTSE'

# Row 18
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range('D18') '0.009499'
$ws.Range('E18').Value = '17OneONEBestin24h'

# Row 19
Set-TextValue $ws.Range('D19') '0.3197'

# Row 20
Set-TextValue $ws.Range('D20') '0.03193'

# Row 21
Set-TextValue $ws.Range('D21') '0.1336'

# Row 22
Set-TextValue $ws.Range('D22') '3.740'

# Row 23
Set-TextValue $ws.Range('D23') '0.04683'

# Row 24
Set-TextValue $ws.Range('D24') '0.1328'

# Row 26
Set-TextValue $ws.Range('D26') '0.004244'

# Row 27
Set-TextValue $ws.Range('D27') '0.00009697'
$ws.Range('E27').Value = '26NitroExNTX'

# Row 28
Set-TextValue $ws.Range('D28') '0.0001965'

# Row 40
Set-TextValue $ws.Range('D40') '0.03627'

# Row 41
Set-TextValue $ws.Range('D41') '0.006288'

# Row 42
Set-TextValue $ws.Range('D42') '0.1049'

# Row 43
Set-TextValue $ws.Range('D43') '0.002729'

# Row 44
Set-TextValue $ws.Range('D44') '0.007390'

# Row 45
Set-TextValue $ws.Range('D45') '0.00005300'

# Row 46
Set-TextValue $ws.Range('D46') '0.00000000749'

# Row 47
Set-TextValue $ws.Range('D47') '0.1999'

# Row 48
Set-TextValue $ws.Range('D48') '0.002293'

# Row 49
Set-TextValue $ws.Range('D49') '0.00002098'

# Row 50
Set-TextValue $ws.Range('D50') '0.0001998'
